$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 09:44"

# Update country statistics rows (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6173653
$ws.Range("C4").Value = 417
$ws.Range("D4").Value = 3425814
$ws.Range("E4").Value = 2560613
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 187226

# Row 28: Ucrania
$ws.Range("B28").Value = 121215
$ws.Range("C28").Value = 2141
$ws.Range("D28").Value = 57114
$ws.Range("E28").Value = 61544
$ws.Range("G28").Value = 30
$ws.Range("H28").Value = 2557

# Row 54: Barein
$ws.Range("E54").Value = 2730
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 190

# Row 60: Armenia
$ws.Range("B60").Value = 43781
$ws.Range("C60").Value = 31
$ws.Range("D60").Value = 37722
$ws.Range("E60").Value = 5180
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 879

# Row 62: Uzbekistan
$ws.Range("B62").Value = 41651
$ws.Range("C62").Value = 227
$ws.Range("E62").Value = 2515
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 317

# Row 72: Australia
$ws.Range("D72").Value = 21350
$ws.Range("E72").Value = 3744

# Row 73: El Salvador
$ws.Range("B73").Value = 25729
$ws.Range("C73").Value = 94
$ws.Range("D73").Value = 14292
$ws.Range("E73").Value = 10720
$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 717

# Row 107: Hungria
$ws.Range("B107").Value = 6139
$ws.Range("C107").Value = 178
$ws.Range("D107").Value = 3761
$ws.Range("E107").Value = 1763
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 615

# Row 111: Hong Kong
$ws.Range("E111").Value = 393
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 89

# Row 153: Georgia
$ws.Range("B153").Value = 1487
$ws.Range("C153").Value = 18
$ws.Range("D153").Value = 1240
$ws.Range("E153").Value = 228

# Row 155: Letonia
$ws.Range("B155").Value = 1396
$ws.Range("C155").Value = 3
$ws.Range("E155").Value = 199

# Row 163: Vietnam
$ws.Range("E163").Value = 312
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 33
